$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.870.59'
$ws.Range("E2").Value = '  +0.66%  '

$ws.Range("D3").Value = '3.543.61'
$ws.Range("E3").Value = '  +4.22%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.06'
$ws.Range("E5").Value = '  +3.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.38'
$ws.Range("E6").Value = '  +3.64%  '

$ws.Range("D7").Value = '3.534.94'
$ws.Range("E7").Value = '  +3.99%  '

$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("E9").Value = '  +3.67%  '

$ws.Range("E10").Value = '  +3.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.92'
$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  +4.44%  '

$ws.Range("D13").Value = '4.150.45'
$ws.Range("E13").Value = '  +4.39%  '

$ws.Range("E14").Value = '  +3.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.12'
$ws.Range("E15").Value = '  +5.02%  '

$ws.Range("D16").Value = '3.552.51'
$ws.Range("E16").Value = '  +5.14%  '

$ws.Range("E17").Value = '  +1.20%  '

$ws.Range("D18").Value = '64.805.11'
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("E19").Value = '  +5.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.42'
$ws.Range("E20").Value = '  +7.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.83'
$ws.Range("E21").Value = '  +3.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.12'
$ws.Range("E22").Value = '  +3.20%  '

$ws.Range("E23").Value = '  +6.80%  '

$ws.Range("D24").Value = '3.692.91'
$ws.Range("E24").Value = '  +4.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.17'
$ws.Range("E25").Value = '  +3.92%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("E27").Value = '  +13.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.64'
$ws.Range("E28").Value = '  +8.44%  '

$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  +5.58%  '

$ws.Range("E31").Value = '  +5.28%  '

$ws.Range("B32").Value = 'RenzoRestakedETH'
$ws.Range("C32").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D32").Value = '3.557.29'
$ws.Range("E32").Value = '  +4.18%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.45'
$ws.Range("E33").Value = '  +24.84%  '

$ws.Range("E34").Value = '  +5.07%  '

$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("E36").Value = '  +3.32%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.96'
$ws.Range("E37").Value = '  +5.66%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '169.90'
$ws.Range("E38").Value = '  -0.60%  '

$ws.Range("E39").Value = '  +6.62%  '

$ws.Range("E40").Value = '  +10.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0806'
$ws.Range("E41").Value = '  +7.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.824'
$ws.Range("E42").Value = '  +4.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.86'
$ws.Range("E43").Value = '  +22.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.58'
$ws.Range("E44").Value = '  +2.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.45'
$ws.Range("E46").Value = '  +5.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").Value = '  +10.35%  '

$ws.Range("E48").Value = '  +4.17%  '

$ws.Range("E49").Value = '  +7.07%  '

$ws.Range("D50").Value = '2.448.00'
$ws.Range("E50").Value = '  +12.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.35'
$ws.Range("E51").Value = '  +16.66%  '
